$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 193; this shifts the existing
# rows 193:259 down to 194:260 (keeping all of their original data),
# and grows the used range from A1:R259 to A1:R260.
$ws.Rows("193:193").Insert()

# Populate the newly inserted row 193 with its data (same shape as every
# other data row in this sheet).
$ws.Range("A193").Value = 3
$ws.Range("B193").Value = "Femacal de La Calera"
$ws.Range("C193").Value = "Coquimbo"
$ws.Range("D193").Value = 44559
$ws.Range("E193").Value = 5
$ws.Range("F193").Value = 100112009
$ws.Range("G193").Value = "Acelga"
$ws.Range("H193").Value = "Sin especificar"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 270
$ws.Range("K193").Value = 2200
$ws.Range("L193").Value = 2500
$ws.Range("M193").Value = 2333
$ws.Range("N193").Value = "$/docena de atados (6 kilos)"
$ws.Range("O193").Value = "Provincia de Quillota"
$ws.Range("P193").Value = 389
$ws.Range("Q193").Value = 6
$ws.Range("R193").Value = "Hortaliza"
